# Apply crypto price/volume updates produced by the scheduled GitHub Actions refresh.
# Values in column D that look like plain decimal numbers are written with a leading
# apostrophe so Excel keeps them as text (matching the original inline-string cells)
# instead of silently re-parsing/truncating them as numbers (e.g. "7.200" -> 7.2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.262.19"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "1.665.26"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").Value = "'219.14"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").Value = "'0.5227"
$ws.Range("E6").Value = "  -1.43%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'0.2665"
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("D9").Value = "'0.06326"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("D10").Value = "'20.93"
$ws.Range("E10").Value = "  -3.61%  "
$ws.Range("D11").Value = "'0.07736"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").Value = "1.677.78"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "'4.441"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").Value = "1.891.84"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "'0.5461"
$ws.Range("E15").Value = "  -1.84%  "
$ws.Range("D16").Value = "0.0₅8218"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").Value = "'64.83"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").Value = "26.275.40"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "'4.656"
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("D21").Value = "'194.78"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "'10.14"
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("D23").Value = "'6.067"
$ws.Range("E23").Value = "  -4.27%  "
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").Value = "'139.61"
$ws.Range("E25").Value = "  -2.06%  "
$ws.Range("D26").Value = "'0.1238"
$ws.Range("E26").Value = "  -3.71%  "
$ws.Range("D27").Value = "'7.200"
$ws.Range("E27").Value = "  -2.91%  "
$ws.Range("D28").Value = "'16.12"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").Value = "'1.415"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("D30").Value = "'0.06153"
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("D31").Value = "'1.284"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").Value = "'3.575"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").Value = "'3.276"
$ws.Range("E33").Value = "  -4.89%  "
$ws.Range("D34").Value = "'1.629"
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("D35").Value = "'0.9744"
$ws.Range("E35").Value = "  -3.29%  "
$ws.Range("D36").Value = "'2.420"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").Value = "'2.788"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").Value = "'0.5736"
$ws.Range("E38").Value = "  -4.80%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01601"
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.016"
$ws.Range("E40").Value = "  -2.60%  "
$ws.Range("D41").Value = "'0.8557"
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").Value = "1.025.40"
$ws.Range("E43").Value = "  -5.40%  "
$ws.Range("D44").Value = "'100.08"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "1.807.20"
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'57.94"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈108"
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("D48").Value = "'1.006"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").Value = "'8.056"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.484"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05187"
$ws.Range("E51").Value = "  -0.40%  "
